$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 25,9
$arr[0,0] = "model_1_8_0"
$arr[0,1] = 0.9835735939480266
$arr[0,2] = 0.990596811261916
$arr[0,3] = 0.985871346448949
$arr[0,4] = 0.9881301023312076
$arr[0,5] = 1.981029460511927
$arr[0,6] = 0.7926666857764346
$arr[0,7] = 2.707764328935551
$arr[0,8] = 1.693889231407512
$arr[1,0] = "model_1_8_1"
$arr[1,1] = 0.9854294411302723
$arr[1,2] = 0.9892375214584468
$arr[1,3] = 0.9845566358410702
$arr[1,4] = 0.9868741238223053
$arr[1,5] = 1.757213737790569
$arr[1,6] = 0.9072516179241674
$arr[1,7] = 2.959729349807819
$arr[1,8] = 1.873123166734724
$arr[2,0] = "model_1_8_2"
$arr[2,1] = 0.9868212871948223
$arr[2,2] = 0.9877153432368719
$arr[2,3] = 0.9832046692930625
$arr[2,4] = 0.9855436525828225
$arr[2,5] = 1.589356688010656
$arr[2,6] = 1.035567660456632
$arr[2,7] = 3.218834492373718
$arr[2,8] = 2.062987558841694
$arr[3,0] = "model_1_8_3"
$arr[3,1] = 0.9878476670560035
$arr[3,2] = 0.9860985898007727
$arr[3,3] = 0.9818513081690468
$arr[3,4] = 0.9841827285134434
$arr[3,5] = 1.465574971167505
$arr[3,6] = 1.17185617104667
$arr[3,7] = 3.478206906208958
$arr[3,8] = 2.257197710454496
$arr[4,0] = "model_1_8_24"
$arr[4,1] = 0.9882632210732005
$arr[4,2] = 0.9642626520537646
$arr[4,3] = 0.9663267656365004
$arr[4,4] = 0.9675425331741824
$arr[4,5] = 1.41545903297039
$arr[4,6] = 3.012574345152819
$arr[4,7] = 6.453494136572458
$arr[4,8] = 4.631830456261419
$arr[5,0] = "model_1_8_23"
$arr[5,1] = 0.9883517355223359
$arr[5,2] = 0.9646711756277411
$arr[5,3] = 0.9665947157266997
$arr[5,4] = 0.9678396422179752
$arr[5,5] = 1.4047841640512
$arr[5,6] = 2.97813676908525
$arr[5,7] = 6.402141352419665
$arr[5,8] = 4.589431623189954
$arr[6,0] = "model_1_8_22"
$arr[6,1] = 0.9884469781607104
$arr[6,2] = 0.965121324402065
$arr[6,3] = 0.9668905474405172
$arr[6,4] = 0.9681673673748372
$arr[6,5] = 1.393297873506592
$arr[6,6] = 2.940190286568686
$arr[6,7] = 6.345445039558128
$arr[6,8] = 4.542663729349606
$arr[7,0] = "model_1_8_21"
$arr[7,1] = 0.9885490380172598
$arr[7,2] = 0.9656169740063627
$arr[7,3] = 0.9672171210982338
$arr[7,4] = 0.9685287683232265
$arr[7,5] = 1.380989424420388
$arr[7,6] = 2.898408191144629
$arr[7,7] = 6.282857016011515
$arr[7,8] = 4.491090144489937
$arr[8,0] = "model_1_8_4"
$arr[8,1] = 0.9885876664040931
$arr[8,2] = 0.9844403398640587
$arr[8,3] = 0.9805226170259109
$arr[8,4] = 0.9828244155917027
$arr[8,5] = 1.376330829467436
$arr[8,6] = 1.311642739001055
$arr[8,7] = 3.732851304456538
$arr[8,8] = 2.451035239236854
$arr[9,0] = "model_1_8_20"
$arr[9,1] = 0.9886577171296853
$arr[9,2] = 0.966161937492369
$arr[9,3] = 0.9675771153741889
$arr[9,4] = 0.9689267238536731
$arr[9,5] = 1.367882691104776
$arr[9,6] = 2.852469051523593
$arr[9,7] = 6.213863912349512
$arr[9,8] = 4.434300051903488
$arr[10,0] = "model_1_8_19"
$arr[10,1] = 0.9887727068105685
$arr[10,2] = 0.9667607495612667
$arr[10,3] = 0.9679736601863327
$arr[10,4] = 0.9693645874476635
$arr[10,5] = 1.354014901354319
$arr[10,6] = 2.80199060306561
$arr[10,7] = 6.137865878052208
$arr[10,8] = 4.371814894290382
$arr[11,0] = "model_1_8_18"
$arr[11,1] = 0.9888933375641306
$arr[11,2] = 0.9674175550338723
$arr[11,3] = 0.9684101011959564
$arr[11,4] = 0.9698458189947601
$arr[11,5] = 1.339466796559262
$arr[11,6] = 2.746623447128225
$arr[11,7] = 6.054221715268155
$arr[11,8] = 4.303140929426847
$arr[12,0] = "model_1_8_17"
$arr[12,1] = 0.9890185624290375
$arr[12,2] = 0.9681366907830371
$arr[12,3] = 0.9688897661911829
$arr[12,4] = 0.9703738709697765
$arr[12,5] = 1.324364640568209
$arr[12,6] = 2.686001995534355
$arr[12,7] = 5.962293651548515
$arr[12,8] = 4.227785473207934
$arr[13,0] = "model_1_8_5"
$arr[13,1] = 0.9891044276561374
$arr[13,2] = 0.9827810497584204
$arr[13,3] = 0.9792370359086598
$arr[13,4] = 0.9814930243933445
$arr[13,5] = 1.31400926861526
$arr[13,6] = 1.451516990748329
$arr[13,7] = 3.979233642212044
$arr[13,8] = 2.641030913725194
$arr[14,0] = "model_1_8_16"
$arr[14,1] = 0.9891468288850973
$arr[14,2] = 0.9689228348507973
$arr[14,3] = 0.9694163441675758
$arr[14,4] = 0.9709525056957768
$arr[14,5] = 1.308895667778542
$arr[14,6] = 2.619731900347333
$arr[14,7] = 5.861374688837207
$arr[14,8] = 4.145211624751986
$arr[15,0] = "model_1_8_15"
$arr[15,1] = 0.989275967609604
$arr[15,2] = 0.969780435075526
$arr[15,3] = 0.96999349562892
$arr[15,4] = 0.9715854644029382
$arr[15,5] = 1.293321499154481
$arr[15,6] = 2.547438219257685
$arr[15,7] = 5.750763289543283
$arr[15,8] = 4.054885493229783
$arr[16,0] = "model_1_8_14"
$arr[16,1] = 0.9894029033261311
$arr[16,2] = 0.9707137947078534
$arr[16,3] = 0.9706250048161145
$arr[16,4] = 0.9722764624387455
$arr[16,5] = 1.278013013948671
$arr[16,6] = 2.468758198362442
$arr[16,7] = 5.629734201788983
$arr[16,8] = 3.95627688139188
$arr[17,0] = "model_1_8_6"
$arr[17,1] = 0.9894484026423797
$arr[17,2] = 0.9811508268275067
$arr[17,3] = 0.9780074481520298
$arr[17,4] = 0.9802061102719669
$arr[17,5] = 1.272525782862559
$arr[17,6] = 1.588940948058739
$arr[17,7] = 4.214884821191592
$arr[17,8] = 2.824679503862495
$arr[18,0] = "model_1_8_13"
$arr[18,1] = 0.9895234955955231
$arr[18,2] = 0.9717264936080252
$arr[18,3] = 0.9713143923333202
$arr[18,4] = 0.9730288612262603
$arr[18,5] = 1.26346955035598
$arr[18,6] = 2.383390063865948
$arr[18,7] = 5.497612699824328
$arr[18,8] = 3.848906098638959
$arr[19,0] = "model_1_8_12"
$arr[19,1] = 0.9896323719624381
$arr[19,2] = 0.97282220010678
$arr[19,3] = 0.9720656099169158
$arr[19,4] = 0.9738462878358919
$arr[19,5] = 1.2503390280902
$arr[19,6] = 2.291024584118203
$arr[19,7] = 5.353641431169495
$arr[19,8] = 3.732255545268034
$arr[20,0] = "model_1_8_7"
$arr[20,1] = 0.9896597237969076
$arr[20,2] = 0.9795717984890543
$arr[20,3] = 0.9768419986200784
$arr[20,4] = 0.9789757425176786
$arr[20,5] = 1.247040388709705
$arr[20,6] = 1.722049321680845
$arr[20,7] = 4.43824387365826
$arr[20,8] = 3.000258666195077
$arr[21,0] = "model_1_8_11"
$arr[21,1] = 0.9897222644479405
$arr[21,2] = 0.9740032627825208
$arr[21,3] = 0.9728817654813128
$arr[21,4] = 0.9747314494956566
$arr[21,5] = 1.239497967574868
$arr[21,6] = 2.191463779485837
$arr[21,7] = 5.197224762295066
$arr[21,8] = 3.605938887334881
$arr[22,0] = "model_1_8_8"
$arr[22,1] = 0.9897702231285891
$arr[22,2] = 0.9780594285028322
$arr[22,3] = 0.9757454984411282
$arr[22,4] = 0.9778097939951902
$arr[22,5] = 1.233714136409889
$arr[22,6] = 1.849538553050958
$arr[22,7] = 4.648388744187128
$arr[22,8] = 3.16664490655931
$arr[23,0] = "model_1_8_10"
$arr[23,1] = 0.9897838405313474
$arr[23,2] = 0.9752706226512532
$arr[23,3] = 0.973765816563899
$arr[23,4] = 0.9756865074444611
$arr[23,5] = 1.232071873582921
$arr[23,6] = 2.084628324533677
$arr[23,7] = 5.027795879511958
$arr[23,8] = 3.469647706063514
$arr[24,0] = "model_1_8_9"
$arr[24,1] = 0.9898049876464278
$arr[24,2] = 0.9766235721686976
$arr[24,3] = 0.9747200384401535
$arr[24,4] = 0.9767126753385619
$arr[24,5] = 1.229521525208087
$arr[24,6] = 1.97057786357973
$arr[24,7] = 4.844918725006354
$arr[24,8] = 3.323208807099505
$ws.Range("A2:I26").Value = $arr
